# refs #882 Video Wall -> HSR Videowall
$wb = $excel.ActiveWorkbook

$wsRisiken = $wb.Worksheets.Item("Risiken")
$wsAenderung = $wb.Worksheets.Item("Änderungsgeschichte")

$wsRisiken.Range("C4").Value = "Die Hardware für die Videowall kann nicht  rechtzeitig geliefert werden."
$wsRisiken.Range("C6").Value = "Kinect erkennt Menschen nicht, die sich parallel zur Wand ausgerichtet vor der Videowall bewegen."
$wsRisiken.Range("B8").Value = "Auflösung der Videowall ungenügend"
$wsRisiken.Range("C8").Value = "Die Auflösung der Videowall ist für das Lesen der Bachelor Posters ungenügend."
$wsRisiken.Range("G8").Value = "Suchen einer Lösung für die Erstellung der Videowall, bei der jeder einzelne Monitor Full HD ist."

$wsRisiken.Range("H8").Select()
$wsAenderung.Select()
